$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 7 new rows for the additional reproduction steps -----------
# Old rows 8-23 shift down to new rows 15-30; this matches the target
# dimension growing from A1:B23 to A1:B30.
$ws.Rows("8:14").Insert()

# Copy the formatting (borders/fill/font) of the existing step rows
# (A6:B7, style ids 5/6) onto the freshly inserted rows so they look the
# same as the other "step" rows in the reproduction-steps block.
$ws.Range("A6:B7").Copy()
$ws.Range("A8:B14").PasteSpecial(-4122)

# --- Title / Description ------------------------------------------------
$ws.Range("B1").Value = "Error al bloquear al usuario"
$ws.Range("B3").Value = "Debería informar que el usuario se encuentra bloquedo"

# --- Reproduction Steps ---------------------------------------------------
$ws.Range("B5").Value = "1- Abrir la aplicación"
$ws.Range("B6").Value = "2- Clic en Sign in"
$ws.Range("B7").Value = "3- Completar campo ""Email address"""
$ws.Range("B8").Value = "4- Completar campo ""Password"""
$ws.Range("B9").Value = "5- Clic en Sign in"
$ws.Range("B10").Value = "6- Completar campo ""Email address"""
$ws.Range("B11").Value = "7- Completar campo ""Password"""
$ws.Range("B12").Value = "8- Clic en Sign in"
$ws.Range("B13").Value = "9- Completar campo ""Email address"""
$ws.Range("B14").Value = "10- Completar campo ""Password"""
$ws.Range("B15").Value = "11- Clic en Sign in"

# --- Expected / Actual behavior ------------------------------------------
$ws.Range("B17").Value = "La página muestra ""Your user has been blockled"", bloquea al usuario"""
$ws.Range("B19").Value = "La página muestra ""User or Password are invalid, Please try again"""

# --- Incidence / severity --------------------------------------------------
$ws.Range("B21").Value = "Incidencia - No se bloquea la cuenta si ingresa de forma erronea 3 veces consecutivas."
$ws.Range("B22").Value = "100% reproduction rate"

# --- Story / Acceptance criteria -------------------------------------------
$ws.Range("B24").Value = "Affected User Story: MT-11, MT-12 y MT-10"

# --- Browsers tested / Environment (unchanged text, left as-is) -----------
$ws.Range("B26").Value = "Google Chrome, v107"
$ws.Range("B28").Value = "QA environment, http://127.0.0.1:5000/"

# --- Notes -------------------------------------------------------------------
# B30 carries a "quotePrefix" style, which a plain .Value assignment would
# otherwise discard. Preserve it by round-tripping the format through a
# scratch cell well outside the used range.
$helper = $ws.Range("Z100")
$ws.Range("B30").Copy()
$helper.PasteSpecial(-4122)
$ws.Range("B30").Value = "Al leer la historia de usuario se entendía que se notificaría que usuario ha sido bloqueado y que por lo tanto no podría intentar iniciar sesión otra vez."
$helper.Copy()
$ws.Range("B30").PasteSpecial(-4122)
$helper.Clear()

# --- Column width / selection tweaks ---------------------------------------
$ws.Columns("B").ColumnWidth = 71.42578125

$ws.Application.Goto($ws.Range("A20"))
$ws.Range("B24").Select()
